# Updates the cryptos list (Price / Volume(1h) columns, plus a few re-ranked
# coin rows) to match the freshly scraped values from the GitHub Actions run.
#
# Numeric-looking "Price" values must stay stored as plain text (as the
# original workbook does), so for those cells we briefly force a Text
# number format before writing the value and then restore the default
# "Normal" style afterwards so no stray style index is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.543.18"
$ws.Range("E2").Value = "  -1.06%  "
$ws.Range("D3").Value = "2.666.17"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.34%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.623"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.52%  "
$ws.Range("E9").Value = "  +2.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.402"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.95%  "
$ws.Range("E11").Value = "  -4.23%  "
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.25"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.97%  "
$ws.Range("E14").Value = "  -5.56%  "
$ws.Range("D15").Value = "3.145.20"
$ws.Range("E15").Value = "  -1.91%  "
$ws.Range("D16").Value = "65.404.18"
$ws.Range("E16").Value = "  -1.02%  "
$ws.Range("D17").Value = "2.667.20"
$ws.Range("E17").Value = "  -2.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("E19").Value = "  -2.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "349.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.59%  "
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.42%  "
$ws.Range("E24").Value = "  +2.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.90%  "
$ws.Range("E26").Value = "  -2.01%  "
$ws.Range("E27").Value = "  -3.20%  "
$ws.Range("E28").Value = "  -9.21%  "
$ws.Range("E29").Value = "  -3.19%  "
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "532.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.99%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.79%  "
$ws.Range("E33").Value = "  -3.61%  "
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.44"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.40%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.47%  "
$ws.Range("E36").Value = "  -3.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.95%  "
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "158.27"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.32%  "
$ws.Range("E40").Value = "  -3.88%  "
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "164.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.13"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.85%  "
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0606"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.641"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.52%  "
$ws.Range("E49").Value = "  -3.58%  "
$ws.Range("E50").Value = "  +1.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.45%  "
